# Daily attendance processing
#
# Normalises the "Recorded By" column (G) on the session-analysis sheet.
# Historically the recorder list was written with the automated "System"
# entry first (e.g. "System, dnasr281@gmail.com"); today's processing pass
# re-orders each comma-separated recorder list so "System" trails the
# human/other recorders instead of leading them
# (e.g. "dnasr281@gmail.com, System").
#
# Only cells whose text starts with "System" AND contains more than one
# recorder (i.e. has a comma) are touched - single-value cells such as
# "System" alone, or lists that don't start with "System" (e.g.
# "admin@admin.com, dnasr281@gmail.com"), are left exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ($text -and $text.StartsWith("System") -and $text.Contains(",")) {
        $parts = $text -split ", "
        $reversed = $parts[($parts.Count - 1)..0]
        $newText = [string]::Join(", ", $reversed)
        $cell.Value = $newText
    }
}
